$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NumberError")

$timestamps = @(
    "Sat Oct 12 23:42:58 EDT 2024",
    "Sat Oct 12 23:43:08 EDT 2024",
    "Sat Oct 12 23:43:18 EDT 2024",
    "Sat Oct 12 23:43:28 EDT 2024",
    "Sat Oct 12 23:43:38 EDT 2024",
    "Sat Oct 12 23:43:47 EDT 2024",
    "Sat Oct 12 23:43:57 EDT 2024",
    "Sat Oct 12 23:44:07 EDT 2024",
    "Sat Oct 12 23:44:16 EDT 2024",
    "Sat Oct 12 23:44:26 EDT 2024",
    "Sat Oct 12 23:44:36 EDT 2024",
    "Sat Oct 12 23:44:46 EDT 2024",
    "Sat Oct 12 23:44:57 EDT 2024",
    "Sat Oct 12 23:45:06 EDT 2024",
    "Sat Oct 12 23:45:16 EDT 2024",
    "Sat Oct 12 23:45:26 EDT 2024",
    "Sat Oct 12 23:45:36 EDT 2024",
    "Sat Oct 12 23:45:45 EDT 2024",
    "Sat Oct 12 23:45:56 EDT 2024",
    "Sat Oct 12 23:46:06 EDT 2024",
    "Sat Oct 12 23:46:16 EDT 2024",
    "Sat Oct 12 23:46:25 EDT 2024",
    "Sat Oct 12 23:46:35 EDT 2024",
    "Sat Oct 12 23:46:45 EDT 2024",
    "Sat Oct 12 23:46:54 EDT 2024",
    "Sat Oct 12 23:47:04 EDT 2024",
    "Sat Oct 12 23:47:14 EDT 2024",
    "Sat Oct 12 23:47:23 EDT 2024",
    "Sat Oct 12 23:47:33 EDT 2024",
    "Sat Oct 12 23:47:43 EDT 2024",
    "Sat Oct 12 23:47:53 EDT 2024",
    "Sat Oct 12 23:48:03 EDT 2024",
    "Sat Oct 12 23:48:13 EDT 2024",
    "Sat Oct 12 23:48:23 EDT 2024",
    "Sat Oct 12 23:48:33 EDT 2024",
    "Sat Oct 12 23:48:42 EDT 2024",
    "Sat Oct 12 23:48:52 EDT 2024",
    "Sat Oct 12 23:49:02 EDT 2024",
    "Sat Oct 12 23:49:11 EDT 2024",
    "Sat Oct 12 23:49:23 EDT 2024",
    "Sat Oct 12 23:49:33 EDT 2024",
    "Sat Oct 12 23:49:43 EDT 2024",
    "Sat Oct 12 23:49:53 EDT 2024",
    "Sat Oct 12 23:50:02 EDT 2024",
    "Sat Oct 12 23:50:12 EDT 2024",
    "Sat Oct 12 23:50:22 EDT 2024",
    "Sat Oct 12 23:50:32 EDT 2024",
    "Sat Oct 12 23:50:42 EDT 2024",
    "Sat Oct 12 23:50:53 EDT 2024",
    "Sat Oct 12 23:51:03 EDT 2024",
    "Sat Oct 12 23:51:13 EDT 2024",
    "Sat Oct 12 23:51:22 EDT 2024",
    "Sat Oct 12 23:51:32 EDT 2024",
    "Sat Oct 12 23:51:42 EDT 2024",
    "Sat Oct 12 23:51:51 EDT 2024",
    "Sat Oct 12 23:52:01 EDT 2024",
    "Sat Oct 12 23:52:11 EDT 2024",
    "Sat Oct 12 23:52:20 EDT 2024",
    "Sat Oct 12 23:52:30 EDT 2024",
    "Sat Oct 12 23:52:40 EDT 2024"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $timestamps[$i]
}
